$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new paragraph at the very start of the document body:
#    "GITHUB: " (bold) + the repo URL (not bold), line spacing 360/auto.
# -----------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$startRange = $firstPara.Range
$startRange.Collapse(1)
$startRange.InsertParagraphBefore() | Out-Null

$newParaRange = $d.Paragraphs.Item(1).Range
$githubOoxml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">GITHUB: </w:t></w:r>
<w:r><w:t>https://github.com/JackTVN/rhythm_game_database_0283473</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$newParaRange.InsertXML($githubOoxml)

# -----------------------------------------------------------------
# 2) Merge the two runs " and their alias " + "(as in, alternative
#    name that they use in a specific " into a single run.
# -----------------------------------------------------------------
$mergedText = " and their alias (as in, alternative name that they use in a specific "
$d.Content.Find.Execute(
    $mergedText, $true, $false, $false, $false, $false, $true, 1, $false,
    $mergedText, 2) | Out-Null

# -----------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the "Album" run to the
#    "Level" run.
# -----------------------------------------------------------------
# 3a) Add it to the run holding "Level".
$levelPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Level`r") {
        $levelPara = $p
        break
    }
}
$levelRange = $levelPara.Range
$levelRunRange = $d.Range($levelRange.Start, $levelRange.End - 1)
$levelOoxml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r w:rsidRPr="004469FF"><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Level</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$levelRunRange.InsertXML($levelOoxml)

# 3b) Remove it from the run holding "Album".
$albumPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Album`r") {
        $albumPara = $p
        break
    }
}
$albumRange = $albumPara.Range
$albumRunRange = $d.Range($albumRange.Start, $albumRange.End - 1)
$albumOoxml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r w:rsidRPr="004469FF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Album</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$albumRunRange.InsertXML($albumOoxml)

Write-Output "done"
